$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '244.38'
$ws.Range('D2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '21.87'
$ws.Range('D3').ClearFormats()

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.399'
$ws.Range('D4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.06043'
$ws.Range('D5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.395'
$ws.Range('D6').ClearFormats()

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8150'
$ws.Range('D7').ClearFormats()

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9262'
$ws.Range('D8').ClearFormats()

$ws.Range('B9').Value = 'One'

$ws.Range('C9').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.01123'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').Value = '8OneONEBestin24h'

$ws.Range('B10').Value = 'WazirX'

$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1439'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').Value = '9WazirXWRX'

$ws.Range('B11').Value = 'MandalaExchangeToken'

$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07456'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'

$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03407'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range('B13').Value = 'BitrueCoin'

$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03049'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').Value = '12BitrueCoinBTR'

$ws.Range('B14').Value = 'BitMartToken'

$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09431'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').Value = '13BitMartTokenBMX'

$ws.Range('B15').Value = 'MCDex'

$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.010'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').Value = '14MCDexMCB'

$ws.Range('B16').Value = 'BitForexToken'

$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001597'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').Value = '15BitForexTokenBF'

$ws.Range('B17').Value = 'CoinExToken'

$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04826'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').Value = '16CoinExTokenCET'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.005688'
$ws.Range('D18').ClearFormats()

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.004155'
$ws.Range('D19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0009910'
$ws.Range('D20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.658'
$ws.Range('D21').ClearFormats()

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1323'
$ws.Range('D25').ClearFormats()

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.00008406'
$ws.Range('D26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0002902'
$ws.Range('D27').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03997'
$ws.Range('D40').ClearFormats()

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006435'
$ws.Range('D41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1074'
$ws.Range('D42').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002902'
$ws.Range('D43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.005822'
$ws.Range('D44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005243'
$ws.Range('D45').ClearFormats()

$ws.Range('E47').Value = '46CoinbaseStockTokenCOIN'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002101'
$ws.Range('D49').ClearFormats()
